$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header cells: "<Name>_old" -> "<Name>_FV2410" and "<Name>_new" -> "<Name>_FV2504"
$oldCols = @("A","B","C","D","E","F","G","H","I","J")
foreach ($col in $oldCols) {
    $cell = $ws.Range($col + "1")
    $cell.Value2 = ($cell.Value2 -replace "_old$", "_FV2410")
}

$newCols = @("L","M","N","O","P","Q","R","S","T","U")
foreach ($col in $newCols) {
    $cell = $ws.Range($col + "1")
    $cell.Value2 = ($cell.Value2 -replace "_new$", "_FV2504")
}

# 2. Freeze the header row (pane split after row 1)
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select()

# 3. Turn the used range into an Excel Table, preserving the existing header
#    formatting (bold / fill / border / wrap / center) instead of inheriting
#    it as an automatic header-style override (dxf).
$headerRange = $ws.Range("A1:U1")
$headerRange.Style = "Normal"

$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U61"), [Type]::Missing, 1)
$lo.Name = "Table1"

$headerRange.Font.Bold = $true
$headerRange.Interior.Color = 14277081
$headerRange.HorizontalAlignment = -4108
$headerRange.WrapText = $true
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2

$lo.TableStyle = ""
